$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update / extend the date column and add new diary notes ---

$ws.Range("A36").Value = 43651
$ws.Range("B36").Value = 'Looking at example dissertations'

$ws.Range("A37").Value = 43654
$ws.Range("B37").Value = 'Checking verified user script works, adding rolling comparison to main - it now compares some stats with the previous run using txt files as storage, identifying buiness value (Probably in tracker as as a way to track the impact of social media campaigns, identify users who are helping spread or should be contacted to help spread, identify the sentiment of the campaign)'

$ws.Range("A38").Value = 43655

$ws.Range("A39").Value = 43656
$ws.Range("A40").Value = 43657
$ws.Range("A41").Value = 43658
$ws.Range("A42").Value = 43661
$ws.Range("A43").Value = 43662
$ws.Range("A44").Value = 43663
$ws.Range("A45").Value = 43664
$ws.Range("A46").Value = 43665
$ws.Range("A47").Value = 43668
$ws.Range("A48").Value = 43669
$ws.Range("A49").Value = 43670
$ws.Range("A50").Value = 43671
$ws.Range("A51").Value = 43672
$ws.Range("A52").Value = 43675
$ws.Range("A53").Value = 43676
$ws.Range("A54").Value = 43677
$ws.Range("A55").Value = 43678
$ws.Range("A56").Value = 43679
$ws.Range("A57").Value = 43682
$ws.Range("A58").Value = 43683
$ws.Range("A59").Value = 43684
$ws.Range("A60").Value = 43685

$ws.Range("A61").Value = 43686
$ws.Range("B61").Value = 'Barrachd placement ends.'

$ws.Range("A62").Value = 43687
$ws.Range("A63").Value = 43688
$ws.Range("A64").Value = 43689
$ws.Range("A65").Value = 43690
$ws.Range("A66").Value = 43691
$ws.Range("A67").Value = 43692
$ws.Range("A68").Value = 43693
$ws.Range("A69").Value = 43694
$ws.Range("A70").Value = 43695
$ws.Range("A71").Value = 43696
$ws.Range("A72").Value = 43697
$ws.Range("A73").Value = 43698
$ws.Range("A74").Value = 43699

# Brand-new rows appended at the bottom of the log
$ws.Range("A75").Value = 43700
$ws.Range("A76").Value = 43701
$ws.Range("A77").Value = 43702
$ws.Range("A78").Value = 43703
$ws.Range("A79").Value = 43704
$ws.Range("A80").Value = 43705
$ws.Range("A81").Value = 43706
$ws.Range("A82").Value = 43707
$ws.Range("A83").Value = 43708
$ws.Range("A84").Value = 43709
$ws.Range("A85").Value = 43710
$ws.Range("A86").Value = 43711
$ws.Range("A87").Value = 43712
$ws.Range("A88").Value = 43713
$ws.Range("A89").Value = 43714

# Give the freshly appended date cells the same date-number-format style
# used by the rest of column A (copy format only, values are untouched)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A75:A89").PasteSpecial(-4122) | Out-Null

# B38 is filled in last so it lands on the next shared-string slot after B61
$ws.Range("B38").Value = 'Modifying verified users script as it failed on a non-existant user over night. Set it running again. Reading example dissertations'

$ws.Range("B38").Select()
